$d = $word.ActiveDocument

# The document currently ends with a paragraph reading:
#   "晴天，今天还是没有课的一天，明天就要上课了。今天就出去玩一下叭~"
# (paragraph 6 of 6). The edit:
#   1. Inserts that same sentence as its own new paragraph right before it
#      (i.e. right after the "2023年3月5日星期天" paragraph).
#   2. Inserts a new date paragraph "2023年3月6日星期一" after that.
#   3. Replaces the text of the (still) final paragraph with the new entry
#      for Monday: "晴天，今天是星期一，又是早八课的一天。"
#
# Anchoring the two new paragraphs off paragraph 5 (instead of splitting
# paragraph 6 itself) makes them inherit paragraph-mark run formatting
# (w:rFonts hint="eastAsia") matching the rest of the diary entries, rather
# than the special "hint=default/eastAsiaTheme" formatting that belongs
# only to the last paragraph.

$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item(6)
$newPara1.Range.InsertBefore("晴天，今天还是没有课的一天，明天就要上课了。今天就出去玩一下叭~")

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item(7)
$newPara2.Range.InsertBefore("2023年3月6日星期一")

# Now update the text of the original last paragraph (still last, now #8).
# Scope the Find/Replace to just that paragraph's range so the identical
# sentence that now also lives in the newly-inserted paragraph #6 is left
# untouched.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute("晴天，今天还是没有课的一天，明天就要上课了。今天就出去玩一下叭~", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "晴天，今天是星期一，又是早八课的一天。", 2)
